$d = $word.ActiveDocument

$replacements = @(
    @("2025-07-07 Monday", "2025-07-08 Tuesday"),
    @("57×22=1254", "49×59=2891"),
    @("33×13=429", "50×63=3150"),
    @("53×42=2226", "45×43=1935"),
    @("35×81=2835", "39×78=3042"),
    @("62×96=5952", "20×39=780"),
    @("25×49=1225", "70×58=4060"),
    @("82×26=2132", "17×45=765"),
    @("63×74=4662", "88×54=4752"),
    @("74×74=5476", "20×37=740"),
    @("25×30=750", "22×66=1452"),
    @("33×47=1551", "47×11=517"),
    @("47×40=1880", "15×63=945"),
    @("64×38=2432", "44×23=1012"),
    @("39×40=1560", "74×85=6290"),
    @("21×38=798", "39×74=2886"),
    @("12×55=660", "21×19=399"),
    @("98×99=9702", "25×78=1950"),
    @("66×27=1782", "78×32=2496"),
    @("44×42=1848", "13×94=1222"),
    @("70×82=5740", "33×49=1617"),
    @("12×84=1008", "63×96=6048"),
    @("51×48=2448", "84×75=6300"),
    @("55×78=4290", "74×12=888"),
    @("91×11=1001", "41×65=2665"),
    @("12×93=1116", "34×92=3128")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
